{"js": "// Resume edit (\"Adding Honest, adaptable, and empathetic to the profile\"):\n//   1) PROFILE paragraph: replace the leading \"A\" with\n//      \"Honest, adaptable, and empathetic\" and swap \"developer\" for\n//      \"software\", i.e.\n//        \"A backend and full stack developer engineer passionate about...\"\n//      becomes\n//        \"Honest, adaptable, and empathetic backend and full stack\n//         software engineer passionate about...\"\n//   2) Section heading \"languages\" -> \"LANGUAGES\" (the displayed text was\n//      already all-caps via paragraph-style formatting; this updates the\n//      underlying run text to match).\n\nconst body = context.document.body;\n\n// 1) PROFILE paragraph rewrite (idempotent: only fires if the old phrase\n// is still present).\nconst profileResults = body.search(\n  \"A backend and full stack developer engineer\",\n  { matchCase: true, matchWholeWord: false }\n);\nprofileResults.load(\"text\");\nawait context.sync();\n\nif (profileResults.items.length > 0) {\n  profileResults.items[0].insertText(\n    \"Honest, adaptable, and empathetic backend and full stack software engineer\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 2) \"languages\" heading -> \"LANGUAGES\" (matchCase so it only matches the\n// stored lowercase run text, not any already-uppercased text).\nconst headingResults = body.search(\"languages\", { matchCase: true, matchWholeWord: true });\nheadingResults.load(\"text\");\nawait context.sync();\n\nif (headingResults.items.length > 0) {\n  headingResults.items[0].insertText(\"LANGUAGES\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Resume edit (\"Adding Honest, adaptable, and empathetic to the profile\"):\n#   1) PROFILE paragraph: replace the leading \"A\" with\n#      \"Honest, adaptable, and empathetic\" and swap \"developer\" for\n#      \"software\", i.e.\n#        \"A backend and full stack developer engineer passionate about...\"\n#      becomes\n#        \"Honest, adaptable, and empathetic backend and full stack\n#         software engineer passionate about...\"\n#   2) Section heading \"languages\" -> \"LANGUAGES\" (the displayed text was\n#      already all-caps via paragraph-style formatting; this updates the\n#      underlying run text to match).\n\n$d = $word.ActiveDocument\n\n# 1) PROFILE paragraph rewrite (idempotent: only fires if the old phrase\n# is still present).\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.Text = \"A backend and full stack developer engineer\"\n$found = $rng.Find.Execute()\nif ($found) {\n    $rng.Text = \"Honest, adaptable, and empathetic backend and full stack software engineer\"\n}\n\n# 2) \"languages\" heading -> \"LANGUAGES\" (MatchCase + MatchWholeWord so it\n# only matches the stored lowercase run text, not any already-uppercased\n# text).\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.MatchCase = $true\n$rng2.Find.MatchWholeWord = $true\n$rng2.Find.Text = \"languages\"\n$found2 = $rng2.Find.Execute()\nif ($found2) {\n    $rng2.Text = \"LANGUAGES\"\n}\n"}
